$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.870.99'
$ws.Range('E2').Value = '  +6.11%  '

$ws.Range('D3').Value = '3.111.32'
$ws.Range('E3').Value = '  +3.89%  '

$ws.Range('E4').Value = '  -0.16%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.33'
$ws.Range('E5').Value = '  +4.73%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.57'
$ws.Range('E6').Value = '  +4.83%  '

$ws.Range('E7').Value = '  -0.13%  '

$ws.Range('D8').Value = '3.099.17'
$ws.Range('E8').Value = '  +3.98%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.533'
$ws.Range('E9').Value = '  +2.05%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.143'
$ws.Range('E10').Value = '  +8.12%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.76'
$ws.Range('E11').Value = '  +12.21%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.470'
$ws.Range('E12').Value = '  +3.09%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000244'
$ws.Range('E13').Value = '  +6.33%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.67'
$ws.Range('E14').Value = '  +6.16%  '

$ws.Range('E15').Value = '  +0.90%  '

$ws.Range('D16').Value = '3.622.53'
$ws.Range('E16').Value = '  +3.79%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.32'
$ws.Range('E17').Value = '  +0.64%  '

$ws.Range('D18').Value = '62.765.51'
$ws.Range('E18').Value = '  +5.91%  '

$ws.Range('D19').Value = '3.099.97'
$ws.Range('E19').Value = '  +3.52%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '453.69'
$ws.Range('E20').Value = '  +5.78%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.15'
$ws.Range('E21').Value = '  +3.73%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.739'
$ws.Range('E22').Value = '  +2.92%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.53'
$ws.Range('E23').Value = '  +6.07%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.81'
$ws.Range('E24').Value = '  +4.08%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.42'
$ws.Range('E25').Value = '  +2.42%  '

$ws.Range('E26').Value = '  +0.18%  '

$ws.Range('E27').Value = '  +4.60%  '

$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.70'
$ws.Range('E28').Value = '  +6.27%  '

$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.30'
$ws.Range('E29').Value = '  +6.53%  '

$ws.Range('E30').Value = '  -0.29%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.88'
$ws.Range('E31').Value = '  +13.63%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.115'
$ws.Range('E32').Value = '  +15.78%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.20'
$ws.Range('E33').Value = '  +5.94%  '

$ws.Range('E34').Value = '  +5.14%  '

$ws.Range('D35').Value = '0.0₃0803'
$ws.Range('E35').Value = '  +6.03%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.10'
$ws.Range('E36').Value = '  +3.16%  '

$ws.Range('E37').Value = '  +6.54%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '50.57'
$ws.Range('E38').Value = '  +3.23%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.01'
$ws.Range('E39').Value = '  +10.89%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.87'
$ws.Range('E40').Value = '  +2.65%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '426.51'
$ws.Range('E41').Value = '  +6.83%  '

$ws.Range('D42').Value = '2.943.53'
$ws.Range('E42').Value = '  +7.04%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0371'
$ws.Range('E43').Value = '  +5.70%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.279'
$ws.Range('E44').Value = '  +11.06%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.110'
$ws.Range('E45').Value = '  +1.70%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.18'
$ws.Range('E46').Value = '  +8.91%  '

$ws.Range('E47').Value = '  -0.02%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.73'
$ws.Range('E48').Value = '  +2.03%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.09'
$ws.Range('E49').Value = '  -0.95%  '

$ws.Range('E50').Value = '  +1.44%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.72'
$ws.Range('E51').Value = '  +5.98%  '
